$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows 4-7 (previously DK1_SmallDecentral, DK2_Central, DK2_LargeDecentral, DK2_SmallDecentral)
$ws.Range("A4:B7").EntireRow.Delete()

# Update first data row: DK1 -> DK, DK1_Central -> DK_Central
$ws.Range("A2").Value = "DK"
$ws.Range("B2").Value = "DK_Central"

# Update second data row: DK1 -> DK, DK1_LargeDecentral -> DK_Decentral
$ws.Range("A3").Value = "DK"
$ws.Range("B3").Value = "DK_Decentral"
